$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix production file errors: correct swapped/incorrect angle values in column D
$ws.Range("D8").Value = 90
$ws.Range("D16").Value = 270
$ws.Range("D17").Value = 270
$ws.Range("D18").Value = 270

# Scroll the view down so row 7 is visible at the top, and leave D8 selected
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D8").Select()
